$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Copy direct-formatting styles from cells that already carry the exact
# resolved style we need, so no new font/cellXf entries are created. ---
# s="3" (bold Calibri 11, row default) lives on any blank cell in row 1
# beyond the old C1 (e.g. D1) before we overwrite its value.
$ws.Cells.Item(1, 4).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null

# s="1" (regular Arial 10, FF222222) lives on the old B3 cell.
$ws.Cells.Item(3, 2).Copy() | Out-Null
$ws.Cells.Item(3, 3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 4).PasteSpecial(-4122) | Out-Null

# The new B3 (Barcode column) must go back to the default (no-style)
# format -- borrow it from a plain cell like A2 before B3's value changes.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 2).PasteSpecial(-4122) | Out-Null

# s="2" (bold Arial 10, FF222222) lives on the old C1 cell; propagate it to
# all the other header cells before C1's own value changes.
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Cells.Item(1, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 3).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Write cell values. Order chosen so the shared-string table comes out
# in the same sequence as the target workbook. ---
$ws.Cells.Item(6, 1).Value = "MatrixTube075"
$ws.Cells.Item(1, 4).Value = "Spacer Length"
$ws.Cells.Item(2, 1).Value = "Eppendorf96"
$ws.Cells.Item(1, 3).Value = "UMI Length"
$ws.Cells.Item(1, 1).Value = "Vessel Type"
$ws.Cells.Item(1, 5).Value = "Location"
$ws.Cells.Item(1, 2).Value = "Barcode"
$ws.Cells.Item(6, 5).Value = "Before First Read"
$ws.Cells.Item(7, 5).Value = "Before Second Read"

$ws.Cells.Item(3, 1).Value = "Eppendorf96"
$ws.Cells.Item(4, 1).Value = "Eppendorf96"
$ws.Cells.Item(5, 1).Value = "Eppendorf96"
$ws.Cells.Item(7, 1).Value = "MatrixTube075"

$ws.Cells.Item(2, 5).Value = "Inline First Read"
$ws.Cells.Item(3, 5).Value = "Before Second Index Read"
$ws.Cells.Item(4, 5).Value = "Inline Second Read"
$ws.Cells.Item(5, 5).Value = "Inline Second Read"

$ws.Cells.Item(2, 2).Value = 12345
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 4).Value = 3

$ws.Cells.Item(3, 2).Value = 34567
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 1

$ws.Cells.Item(4, 2).Value = 66789
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 1

$ws.Cells.Item(5, 2).Value = 77891
$ws.Cells.Item(5, 3).Value = 9
$ws.Cells.Item(5, 4).Value = 3

$ws.Cells.Item(6, 2).Value = 87654
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 4).Value = 2

$ws.Cells.Item(7, 2).Value = 87654
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = 2

# --- Column widths (character-width units; the XML "width" the engine
# writes is ColumnWidth + 5/6, quantised to the nearest 1/6 character). ---
$ws.Columns.Item(1).ColumnWidth = 11.833333333333334
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 23.666666666666668
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668
$ws.Columns.Item(5).ColumnWidth = 20.0

# --- Selection matches the author's last-saved cursor position. ---
$ws.Range("E7").Select() | Out-Null
